# Active_Outages.xlsx refresh — 6/18/2025, 5:03:20 PM
# - Elapsed Duration(Hrs) values tick forward on all still-open outages
# - A new outage row is appended on the R1 sheet (Haj Removal / HAJ0125 / R5)

$wb = $excel.ActiveWorkbook

# --- R1 sheet: update elapsed durations, then append new outage row ---
$ws = $wb.Worksheets.Item("R1")
$ws.Range("G2").Value = "3930:17:18"
$ws.Range("G3").Value = "69:49:56"
$ws.Range("G4").Value = "92:49:56"

$ws.Range("B6").Value = "R5"
$ws.Range("D6").Value = "HAJ0125"
$ws.Range("I6").Value = "SCECO"
$ws.Range("J6").Value = "Haj Removal"
$ws.Range("L6").Value = "Latis"

# --- R2 sheet: update elapsed durations ---
$ws = $wb.Worksheets.Item("R2")
$ws.Range("G2").Value = "12111:40:58"
$ws.Range("G3").Value = "3241:24:27"
$ws.Range("G4").Value = "479:36:01"

# --- R4 sheet: update elapsed durations ---
$ws = $wb.Worksheets.Item("R4")
$ws.Range("G2").Value = "2957:30:47"
$ws.Range("G3").Value = "184:43:02"
$ws.Range("G4").Value = "72:55:27"
$ws.Range("G5").Value = "70:33:00"

# --- R5 sheet: update elapsed duration ---
$ws = $wb.Worksheets.Item("R5")
$ws.Range("G2").Value = "431:29:46"

# --- R6 sheet: update elapsed duration ---
$ws = $wb.Worksheets.Item("R6")
$ws.Range("G2").Value = "72:02:04"
